# Regenerate merged AHB files
# 1. Rename header-row labels: "_old" -> "_FV2404", "_new" -> "_FV2410"
# 2. Freeze the header row (row 1)
# 3. Wrap the data range in an Excel table (Table1) with AutoFilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$newHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J: "<Label>_old" -> "<Label>_FV2404"
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($oldHeaders[$i])_FV2404"
}

# Column K: "diff" stays the same

# Columns L-U: "<Label>_new" -> "<Label>_FV2410"
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($newHeaders[$i])_FV2410"
}

# Freeze panes at row 1 (split happens below row 1, so row 1 stays visible)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a real Excel table with header row + autofilter
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U70"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
